# Auto-generated edit script applying the crypto price/volume update
# described by the commit "Updated cryptos list ... with GitHub Actions".
# Every D/E cell in this sheet is stored as plain text (t="inlineStr" in
# the original workbook) -- the "Price" column holds dotted-grouping
# strings like "66.712.95" and the "Volume(1h)" column holds
# space-padded percent strings like "  +2.39%  ". We assign plain
# strings to Range.Value to match. For Price cells whose new value would
# otherwise be auto-recognised as a plain number (e.g. "578.87",
# "1.00", "0.0360") we lead with an apostrophe, exactly as typing the
# same text into Excel's UI would, so the cell keeps its original Text
# type instead of silently becoming a Number (which would also perform
# ugly things like dropping the significant trailing zero in "1.00" or
# "0.0360").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.715.05"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "3.082.24"
$ws.Range("E3").Value = "  +4.39%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'578.87"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").Value = "'167.79"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.079.14"
$ws.Range("E8").Value = "  +4.50%  "
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "'6.57"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").Value = "'0.153"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("E12").Value = "  +5.43%  "
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "'36.43"
$ws.Range("E14").Value = "  +6.56%  "
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "3.596.13"
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("D17").Value = "66.768.92"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("D18").Value = "'7.21"
$ws.Range("E18").Value = "  +3.87%  "
$ws.Range("D19").Value = "3.082.91"
$ws.Range("E19").Value = "  +4.51%  "
$ws.Range("E20").Value = "  +10.88%  "
$ws.Range("D21").Value = "'464.55"
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("E23").Value = "  +4.10%  "
$ws.Range("D24").Value = "'83.20"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  +5.70%  "
$ws.Range("D26").Value = "'12.87"
$ws.Range("E26").Value = "  +6.49%  "
$ws.Range("D27").Value = "'10.14"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'7.97"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("E31").Value = "  +3.38%  "
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "'28.15"
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("E34").Value = "  +3.64%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  +2.65%  "
$ws.Range("D37").Value = "'5.88"
$ws.Range("E37").Value = "  +3.12%  "
$ws.Range("E38").Value = "  +7.74%  "
$ws.Range("D39").Value = "'47.10"
$ws.Range("E39").Value = "  +6.46%  "
$ws.Range("E40").Value = "  +6.26%  "
$ws.Range("D41").Value = "'50.18"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("D43").Value = "'8.67"
$ws.Range("E43").Value = "  +2.84%  "
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "'0.0360"
$ws.Range("E45").Value = "  +2.42%  "
$ws.Range("D46").Value = "'382.93"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").Value = "2.759.74"
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("D48").Value = "'134.83"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'24.61"
$ws.Range("E50").Value = "  +5.83%  "
$ws.Range("D51").Value = "'2.23"
$ws.Range("E51").Value = "  +3.15%  "
